# Auto commit at 2026-02-22 10:32:39.89
#
# Appends 8 more days (2026-02-14 .. 2026-02-21) of per-station daily
# charging-station data to the existing 日期/站点/... table -- two rows
# per day (四方坪站 then 高岭站), continuing the exact pattern already
# present in rows 2-27. Also adds the trailing "spacer" row that only
# carries a number format on F, and moves the selection to where the
# user ended up after typing the new data in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: date serial, station name, 充电服务费收入, 充电总收入, 充电量, 总订单数量
$rows = @(
    @(46067, "四方坪站", 8774.75,            8003.16,  3269.31, 365),
    @(46067, "高岭站",   2893.05,            2581.69,  853.54,  94),
    @(46068, "四方坪站", 6862.61,            6336.94,  2540.74, 291),
    @(46068, "高岭站",   1356.01,            1207.93,  425.14,  49),
    @(46069, "四方坪站", 4827.05,            4651.75,  1770.27, 190),
    @(46069, "高岭站",   677.07,             642.65,   162.37,  17),
    @(46070, "四方坪站", 3649.8,             3480.67,  1328.93, 140),
    @(46070, "高岭站",   920.35,             882.33,   247.62,  31),
    @(46071, "四方坪站", 4725.02,            4562.08,  1714.89, 183),
    @(46071, "高岭站",   1205.4,             1119.72,  348.95,  33),
    @(46072, "四方坪站", 5022.96,            4761.92,  1877.18, 189),
    @(46072, "高岭站",   1066.29,            971.52,   273.01,  31),
    @(46073, "四方坪站", 5603.02,            5290.31,  2047.6,  217),
    @(46073, "高岭站",   1219.54,            1140.86,  312.93,  38),
    @(46074, "四方坪站", 5977.755,           5468.88,  2214.36, 255),
    @(46074, "高岭站",   1502.22,            1351.01,  418.56,  56)
)

$startRow = 28
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

# Trailing spacer row below the data (row 44): only F44 is touched, carrying
# the same 2-decimal number format used by columns C/D/E, with no value.
$ws.Range("F44").NumberFormat = "0.00_);[Red]\(0.00\)"

# Leave the selection where the user would have ended up after entering the
# last row of new data.
$ws.Range("E49").Select() | Out-Null
